$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the first student row (Muhammad Akmalul Iman Liari); this shifts the
# second student row (Syifa Khaista Khairunnisa) up so it becomes row 2.
$ws.Rows(2).Delete()

# Add the new "user_id" column (J). Copy the header formatting from the
# existing header row (A1) onto J1 so it matches the other header cells.
$ws.Range("A1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("J1").Value = "user_id"
$ws.Range("J2").Value = "U2020-0056"
